# Add a new note row (row 15) to Plan1 documenting the plan to build the
# basic CRUD first, then lock it down with Spring Security roles.
#
# "Spring Security" is rendered in bold, matching the rich-text run split
# recorded in the target sharedStrings.xml entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "Primeiro vamos focar no básico (CRUD normal). Depois, quando você implementar o "
$bold   = "Spring Security"
$suffix = ", aí a gente bloqueia o acesso via roles"

$cell = $ws.Range("A15")
$cell.Value = $prefix + $bold + $suffix

# 1-based character offsets into the cell text for Range.Characters(Start, Length)
$boldStart = $prefix.Length + 1

$boldRun = $cell.Characters($boldStart, $bold.Length)
$boldRun.Font.Bold = $true

# Touch the tail run's font so Excel materialises explicit run properties
# (rPr) for it too, same as the recorded diff (non-bold but still styled).
$suffixStart = $boldStart + $bold.Length
$suffixRun = $cell.Characters($suffixStart, $suffix.Length)
$suffixRun.Font.Size = 11
$suffixRun.Font.Name = "Calibri"

# Move the active selection onto the newly added cell.
$cell.Select() | Out-Null
